$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (row, Cluster Name, Active cases) data for rows 2..46, sorted alphabetically by
# Cluster Name as in the edited workbook. Row 1 (header: "Cluster Name" / "Active cases")
# is unchanged.
$data = @(
  @(2, '3528 Ottoman Village Aged Care Broadmeadows Outbreak', 24),
  @(3, '3652 Regis Aged Care Dandenong North Outbreak', 25),
  @(4, '3824 Estia Health South Morang Outbreak', 56),
  @(5, '3961 Heritage Water Gardens Aged Care Facility Sydenham Outbreak', 15),
  @(6, 'Aintree Primary School Aintree', 13),
  @(7, 'Alfred Health The Alfred Hospital Melbourne Outbreak', 11),
  @(8, 'Armstrong Creek School Armstrong Creek Outbreak', 12),
  @(9, 'Australian Meat Group Abattoir Dandenong South', 10),
  @(10, 'Berwick Fields Primary School Berwick Outbreak', 13),
  @(11, 'Berwick Lodge Primary School Berwick Outbreak', 21),
  @(12, "Bubup Womindjeka Family and Children's Centre Port Melbourne Outbreak", 11),
  @(13, "CREST Children's Sanctuary Dandenong Outbreak", 11),
  @(14, 'Clifton Hill Primary School Clifton Hill Outbreak', 13),
  @(15, 'Dandenong North Primary School Dandenong Outbreak', 10),
  @(16, 'Elements Childcare Warralily Armstrong Creek Outbreak', 21),
  @(17, 'G & K OConnor PTY LTD Pakenham', 10),
  @(18, 'Hamlyn Views School Hamlyn Heights Outbreak', 10),
  @(19, 'KingKids Early Learning Centre and Kindergarten Hallam Outbreak', 11),
  @(20, 'Lilydale Motor Inn Lilydale Outbreak', 12),
  @(21, 'Lowanna College Newborough Outbreak', 35),
  @(22, 'McQuinns Gym Bendigo Outbreak', 18),
  @(23, 'Mercy Health Werribee Mercy Hospital Outbreak', 10),
  @(24, 'Metcash Limited Distribution Centre Laverton North Outbreak', 15),
  @(25, 'Monash Health Dandenong Hospital Dandenong Outbreak', 11),
  @(26, 'Monash Health Dandenong Hospital Emergency Department Placeholder', 69),
  @(27, 'Northern Bay College Wexford Campus Corio Outbreak', 32),
  @(28, 'Rosewood Downs Special Accommodation Home Dandenong Outbreak', 13),
  @(29, 'Saint Augustines Primary School Wodonga Outbreak', 14),
  @(30, "St Mary's Primary School Swan Hill Outbreak", 20),
  @(31, 'St Thereses Primary School Kennington Outbreak', 14),
  @(32, 'St Vincents Hospital Emergency Department Melbourne Outbreak', 19),
  @(33, 'St. Brendans Catholic Primary School Lakes Entrance Outbreak', 13),
  @(34, 'TUROSI PTY LTD Thomastown Outbreak', 15),
  @(35, "The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1A", 10),
  @(36, "The Royal Children's Hospital Parkville Outbreak", 10),
  @(37, 'Vizzarri Farms Koo Wee Rup Outbreak', 24),
  @(38, 'Werribee Mercy Hospital Emergency Department Outbreak', 37),
  @(39, 'Western Health Sunshine Hospital Emergency Department St Albans Outbreak', 17),
  @(40, 'Wodonga Cemetery Wodonga Outbreak', 41),
  @(41, 'Wodonga Primary School Wodonga Outbreak', 13),
  @(42, 'Wodonga Senior Secondary College Wodonga Outbreak', 19),
  @(43, 'Wodonga South Primary School Wodonga Outbreak', 32),
  @(44, 'Woodend Primary School Woodend Outbreak', 20),
  @(45, 'Yallourn Power Station Yallourn', 10),
  @(46, 'Yooralla Disability Residential Care Alfrieda Street St Albans Outbreak', 13)
)

foreach ($row in $data) {
  $r = $row[0]
  $name = $row[1]
  $cases = $row[2]
  $ws.Cells.Item($r, 1).Value = $name
  $ws.Cells.Item($r, 2).Value = $cases
}
